# Generate Report for Handoff
#
# The only substantive data change introduced by this report regeneration is
# an updated "Latest Handoff Datetime" value for the
# 6ffffe0b-afb0-4436-ab56-02431434309c.md entry on the "zh-cn" localization
# status sheet: it moves from 2016-09-03 18:53:42 to 2016-09-03 18:53:58.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

# Row 5 corresponds to 6ffffe0b-afb0-4436-ab56-02431434309c.md (column A),
# column H is "Latest Handoff Datetime".
$ws.Range("H5").Value = "2016-09-03 18:53:58"
